$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.675.33"
$ws.Range("E2").Value = "  -0.03%  "
$ws.Range("D3").Value = "1.643.32"
$ws.Range("E4").Value = "  +0.29%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.00"
$ws.Range("E5").Value = "  +0.77%  "
$ws.Range("E6").Value = "  +0.89%  "
$ws.Range("E7").Value = "  +0.24%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.253"
$ws.Range("E8").Value = "  +0.21%  "
$ws.Range("E9").Value = "  +0.98%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.28"
$ws.Range("E10").Value = "  +0.26%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0842"
$ws.Range("E11").Value = "  +0.07%  "
$ws.Range("D12").Value = "1.871.79"
$ws.Range("E12").Value = "  +0.67%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.654.08"
$ws.Range("E13").Value = "  +1.28%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.21"
$ws.Range("E14").Value = "  +2.77%  "
$ws.Range("E15").Value = "  +1.21%  "
$ws.Range("E16").Value = "  +3.21%  "
$ws.Range("D17").Value = "26.708.77"
$ws.Range("E17").Value = "  +0.09%  "
$ws.Range("E18").Value = "  +1.21%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "216.61"
$ws.Range("E19").Value = "  -1.04%  "
$ws.Range("E20").Value = "  +0.25%  "
$ws.Range("E21").Value = "  +1.61%  "
$ws.Range("E22").Value = "  +2.34%  "
$ws.Range("E24").Value = "  +12.45%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.81"
$ws.Range("E25").Value = "  -1.21%  "
$ws.Range("E26").Value = "  +0.27%  "
$ws.Range("E27").Value = "  -1.10%  "
$ws.Range("E28").Value = "  +4.57%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.78"
$ws.Range("E29").Value = "  +1.41%  "
$ws.Range("E30").Value = "  +2.32%  "
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("E32").Value = "  +2.65%  "
$ws.Range("E33").Value = "  +2.12%  "
$ws.Range("D34").Value = "1.272.97"
$ws.Range("E34").Value = "  +4.63%  "
$ws.Range("E35").Value = "  +2.36%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0182"
$ws.Range("E36").Value = "  +6.36%  "
$ws.Range("E37").Value = "  +0.23%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.533"
$ws.Range("E38").Value = "  +6.35%  "
$ws.Range("E39").Value = "  +2.95%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.815"
$ws.Range("E41").Value = "  +2.62%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.25"
$ws.Range("E42").Value = "  -1.37%  "
$ws.Range("E43").Value = "  +2.18%  "
$ws.Range("D44").Value = "1.782.40"
$ws.Range("E44").Value = "  +0.64%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.32"
$ws.Range("E45").Value = "  -0.56%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "60.03"
$ws.Range("E46").Value = "  +9.05%  "
$ws.Range("E47").Value = "  +2.96%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0516"
$ws.Range("E48").Value = "  +0.74%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.81"
$ws.Range("E49").Value = "  +2.65%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0972"
$ws.Range("E50").Value = "  +3.18%  "
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.406"
$ws.Range("E51").Value = "  -0.40%  "
